# B6-PowerPoint.pptx edit — Wed, Apr 08, 2020 12:08:47 PM
#
# 1) Three tables (one each on slides 14, 15, 16) get switched from the
#    deck's custom "Table_0" style to a built-in PowerPoint table style
#    (referenced purely by its GUID, no local <a:tblStyle> definition
#    needed).
# 2) The presentation's design swaps back from the "Integral" theme to
#    the plain "Office Theme" colour palette (fonts/effects were already
#    identical between the two themes, only the 12 theme colours and the
#    theme/colour-scheme names differ).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newTableStyleId = "{BFDDFFD1-87AF-4CDC-A6DD-B879B0EE9C97}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the master's colour scheme back to "Office Theme" ------
# RGB(r,g,b) -> the little-endian OLE colour PowerPoint's ColorScheme
# setter expects.
function RGBColor($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$officeThemeColors = @(
    (RGBColor 0x00 0x00 0x00),  # dk1
    (RGBColor 0xFF 0xFF 0xFF),  # lt1
    (RGBColor 0x44 0x54 0x6A),  # dk2
    (RGBColor 0xE7 0xE6 0xE6),  # lt2
    (RGBColor 0x5B 0x9B 0xD5),  # accent1
    (RGBColor 0xED 0x7D 0x31),  # accent2
    (RGBColor 0xA5 0xA5 0xA5),  # accent3
    (RGBColor 0xFF 0xC0 0x00),  # accent4
    (RGBColor 0x44 0x72 0xC4),  # accent5
    (RGBColor 0x70 0xAD 0x47),  # accent6
    (RGBColor 0x05 0x63 0xC1),  # hlink
    (RGBColor 0x95 0x4F 0x72)   # folHlink
)

$master = $p.SlideMaster
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $master.ColorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
